# Applies the edit described by the diff:
#  - Rename "Sheet1" to "intraday"
#  - Add a new reflective note row (row 79, col B) on the intraday sheet,
#    highlighted with an orange fill
#  - Widen columns B and C on the intraday sheet
#  - Update the saved view/selection state on both the intraday and
#    swing_strategy sheets

$wb = $excel.ActiveWorkbook

# --- Rename Sheet1 -> intraday -------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Name = "intraday"

$ws2 = $wb.Worksheets.Item("swing_strategy")

# --- Touch swing_strategy's selection first (diff shows its saved
#     selection moved from A80:A82 to A66) so its view state is captured,
#     then hop back to intraday so that sheet keeps "tabSelected" ---------
$ws2.Activate() | Out-Null
$ws2.Range("A66").Select() | Out-Null

# --- Back to intraday: new note row, widened columns, view state --------
$ws1.Activate() | Out-Null

$ws1.Columns.Item(2).ColumnWidth = 122.4167
$ws1.Columns.Item(3).ColumnWidth = 17.25

$noteCell = $ws1.Range("B79")
$noteCell.Value = "This 2024 year is a bad year because of not managing my risk to reward and I should have control in my trading that is maximum 2 trade per day"
$noteCell.Interior.Color = 49407

$excel.ActiveWindow.Zoom = 68
$ws1.Range("C79").Select() | Out-Null
